$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 12-23 (they will be rebuilt with new layout/content)
$ws.Rows("12:23").Delete()
# Insert 13 fresh rows to host the restructured content (rows 12-24)
$ws.Rows("12:24").Insert()

$ws.Range("A12").Value = 'Docentes responsáveis:'
$ws.Range("B13").Value = '5840671 - Francisco José Moreira Chaves'
$ws.Range("C13").Value = '5840671 - Francisco José Moreira Chaves'
$ws.Range("A14").Value = 'Programa resumido:'
$ws.Range("B14").Value = 'Mercado - Estimativa de investimento: - Análise Econômica de Investimentos'
$ws.Range("C14").Value = 'Mercado - Estimativa de investimento: - Análise Econômica de Investimentos'
$ws.Range("A15").Value = 'Short syllabus:'
$ws.Range("B15").Value = 'Market - Estimated investment : - Economic Analysis of Investments'
$ws.Range("C15").Value = 'Market - Estimated investment : - Economic Analysis of Investments'
$ws.Range("A16").Value = 'Programa:'
$ws.Range("B16").Value = 'Mercado - evolução do mercado - Marketing e análise de mercado - Estimativa de investimento: - capital de giro - capital fixo - bens tangíveis/bens intangíveis - investimentos preliminares/investimentos reais - custos fixos/custos variáveis, depreciação, Conceito econômico de externalidades e abordagens teóricas, Elementos para internalizar as externalidades, Controle direto ou regularização na empresa, métodos indiretos c dados observados, métodos indiretos c dados supostos, métodos diretos c dados supostos, métodos diteros c dados observados, Análise Econômica de Investimentos (aspectos básicos) - Técnicas Estatísticas e taxas de juros - Aspectos básicos de Engenharia Econômica  Distribuição Beta para análise em ambiente de risco.'
$ws.Range("C16").Value = 'Mercado - evolução do mercado - Marketing e análise de mercado - Estimativa de investimento: - capital de giro - capital fixo - bens tangíveis/bens intangíveis - investimentos preliminares/investimentos reais - custos fixos/custos variáveis, depreciação, Conceito econômico de externalidades e abordagens teóricas, Elementos para internalizar as externalidades, Controle direto ou regularização na empresa, métodos indiretos c dados observados, métodos indiretos c dados supostos, métodos diretos c dados supostos, métodos diteros c dados observados, Análise Econômica de Investimentos (aspectos básicos) - Técnicas Estatísticas e taxas de juros - Aspectos básicos de Engenharia Econômica  Distribuição Beta para análise em ambiente de risco.'
$ws.Range("A17").Value = 'Syllabus:'
$ws.Range("B17").Value = 'Market - market trends - Marketing and market analysis - Estimate of investment: - working capital - capital assets - tangible / intangibles assets - Preliminary investments / real investments - fixed costs / variable costs , depreciation , economic concept of externalities and approaches theoretical , elements to internalize externalities , direct control or stabilize the company , indirect methods and observed data , indirect methods and data assumptions , methods and alleged direct data  methods and observed data , Economic Analysis of Investments ( basics aspects) - Statistical Techniques and interest rates - basics of Engineering Economy - Beta Distribution for analysis in the risk environment .'
$ws.Range("C17").Value = 'Market - market trends - Marketing and market analysis - Estimate of investment: - working capital - capital assets - tangible / intangibles assets - Preliminary investments / real investments - fixed costs / variable costs , depreciation , economic concept of externalities and approaches theoretical , elements to internalize externalities , direct control or stabilize the company , indirect methods and observed data , indirect methods and data assumptions , methods and alleged direct data  methods and observed data , Economic Analysis of Investments ( basics aspects) - Statistical Techniques and interest rates - basics of Engineering Economy - Beta Distribution for analysis in the risk environment .'
$ws.Range("A18").Value = 'Avaliação:'
$ws.Range("A19").Value = 'Método:'
$ws.Range("B19").Value = 'Por meio de aulas presenciais, com apresentação dos fundamentos, e resolução de exercícios e exemplos aplicativos com uso de tabelas e normas específicas.'
$ws.Range("C19").Value = 'Por meio de aulas presenciais, com apresentação dos fundamentos, e resolução de exercícios e exemplos aplicativos com uso de tabelas e normas específicas.'
$ws.Range("A20").Value = 'Critério:'
$ws.Range("B20").Value = 'A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2'
$ws.Range("C20").Value = 'A Avaliação será: (P1 + 2P2)/3 Onde: P1: Prova Individual - c/ peso-1 P2: : Prova Individual - c/ peso-2'
$ws.Range("A21").Value = 'Norma de recuperação:'
$ws.Range("B21").Value = 'Prova de exame.'
$ws.Range("C21").Value = 'Prova de exame.'
$ws.Range("A22").Value = 'Bibliografia:'
$ws.Range("B22").Value = 'Engenharia econômica e análise de custos. Henrique Hirschfeld. 7 ed. editora atlas. 2007'
$ws.Range("C22").Value = 'Engenharia econômica e análise de custos. Henrique Hirschfeld. 7 ed. editora atlas. 2007'
$ws.Range("A23").Value = 'Requisitos:'
$ws.Range("B24").Value = 'LOB1012 -  Estatística  (Requisito)
'
$ws.Range("C24").Value = 'LOB1012 -  Estatística  (Requisito)
'

# Clear stray empty cells in rows with fewer populated columns
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("A13").ClearContents()
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("B23").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("A24").ClearContents()

# Set explicit custom row heights where required; AutoFit (back to default) elsewhere
$ws.Rows.Item(12).AutoFit()
$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).AutoFit()
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 120
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).RowHeight = 30
